$wb = $excel.ActiveWorkbook
$eventSheet = $wb.Worksheets.Item("event")
$ws = $wb.Worksheets.Add($null, $eventSheet)
$ws.Name = "milestone"

$ws.Range("A1").Value = "In Portuguese"
$ws.Range("B1").Value = "Name"
$ws.Range("D1").Value = "Query"

$data = @(
    @("Primeiro Dente", "Surgery"),
    @("Balbuciar", "Hospitalization"),
    @("Chamar os Pais", "Fracture"),
    @("Desmaio", "Fainting"),
    @("Intoxicação", "Intoxication"),
    @("Coma", "Coma"),
    @("Reação Alérgica", "Allergic Reaction"),
    @("Convulsão", "Convulsion"),
    @("Contusão", "Contusion"),
    @("Concussão", "Concussion"),
    @("Torção", "Torsion"),
    @("Corte", "Cut"),
    @("Ingestão", "Ingestion"),
    @("Choque Anafilático", "Anaphylactic Shock "),
    @("Queda", "Fall"),
    @("Intubação ", "Intubation")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
}

$ws.Range("D2").Formula = "=""(:Milestone {uuid: apoc.create.uuid(), name: '"" & B2 &""', inPortuguese: '"" & A2 &""'}),"""

$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 30.5546875

$ws.Range("D12").Select()

$eventWs = $wb.Worksheets.Item("event")
$eventWs.Range("B7").Select()
